# Updates hourly crypto snapshot values (price + 1h volume change),
# and fixes the row order for three coin pairs that swapped ranking
# position (rows 39-41 and 45-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.635.20"
$ws.Range("E2").Value = "'  -0.30%  "
$ws.Range("D3").Value = "'3.778.72"
$ws.Range("E3").Value = "'  -1.58%  "
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("D5").Value = "'596.98"
$ws.Range("E5").Value = "'  -0.21%  "
$ws.Range("D6").Value = "'169.26"
$ws.Range("E6").Value = "'  +0.96%  "
$ws.Range("D7").Value = "'3.775.57"
$ws.Range("E7").Value = "'  -1.65%  "
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "'  +0.52%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "'  +2.02%  "
$ws.Range("D11").Value = "'6.50"
$ws.Range("E11").Value = "'  +1.43%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "'  +1.15%  "
$ws.Range("E13").Value = "'  +6.21%  "
$ws.Range("D14").Value = "'36.90"
$ws.Range("E14").Value = "'  +0.32%  "
$ws.Range("D15").Value = "'4.421.48"
$ws.Range("E15").Value = "'  -1.71%  "
$ws.Range("D16").Value = "'3.787.29"
$ws.Range("E16").Value = "'  -1.74%  "
$ws.Range("D17").Value = "'19.01"
$ws.Range("E17").Value = "'  +6.08%  "
$ws.Range("D18").Value = "'67.782.74"
$ws.Range("E18").Value = "'  -0.34%  "
$ws.Range("D19").Value = "'7.26"
$ws.Range("E19").Value = "'  -0.50%  "
$ws.Range("E20").Value = "'  +1.11%  "
$ws.Range("D21").Value = "'10.54"
$ws.Range("E21").Value = "'  -1.57%  "
$ws.Range("D22").Value = "'467.07"
$ws.Range("E22").Value = "'  +0.58%  "
$ws.Range("D23").Value = "'0.725"
$ws.Range("E23").Value = "'  -0.41%  "
$ws.Range("E24").Value = "'  -5.21%  "
$ws.Range("D25").Value = "'83.49"
$ws.Range("E25").Value = "'  +0.80%  "
$ws.Range("E26").Value = "'  +1.55%  "
$ws.Range("D27").Value = "'12.14"
$ws.Range("E27").Value = "'  +1.45%  "
$ws.Range("D28").Value = "'10.36"
$ws.Range("E28").Value = "'  +4.59%  "
$ws.Range("E29").Value = "'  +0.08%  "
$ws.Range("E30").Value = "'  -0.96%  "
$ws.Range("D31").Value = "'3.935.70"
$ws.Range("E31").Value = "'  -1.66%  "
$ws.Range("D32").Value = "'7.62"
$ws.Range("E32").Value = "'  -0.08%  "
$ws.Range("E33").Value = "'  -1.79%  "
$ws.Range("D34").Value = "'30.44"
$ws.Range("E34").Value = "'  -1.82%  "
$ws.Range("D35").Value = "'9.15"
$ws.Range("E35").Value = "'  -2.65%  "
$ws.Range("D36").Value = "'3.753.61"
$ws.Range("E36").Value = "'  -1.63%  "
$ws.Range("D37").Value = "'3.86"
$ws.Range("E37").Value = "'  +6.85%  "
$ws.Range("E38").Value = "'  +1.41%  "
$ws.Range("B39").Value = "'Mantle"
$ws.Range("C39").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "'1.01"
$ws.Range("E39").Value = "'  -0.89%  "
$ws.Range("B40").Value = "'Filecoin"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'5.90"
$ws.Range("E40").Value = "'  +0.59%  "
$ws.Range("B41").Value = "'Kaspa"
$ws.Range("C41").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.138"
$ws.Range("E41").Value = "'  -1.28%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "'  -0.01%  "
$ws.Range("E43").Value = "'  +1.75%  "
$ws.Range("E44").Value = "'  +0.02%  "
$ws.Range("B45").Value = "'Cosmos"
$ws.Range("C45").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").Value = "'8.72"
$ws.Range("E45").Value = "'  +1.92%  "
$ws.Range("B46").Value = "'Stacks"
$ws.Range("C46").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.96"
$ws.Range("E46").Value = "'  +0.12%  "
$ws.Range("B47").Value = "'OKB"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'46.39"
$ws.Range("E47").Value = "'  -1.26%  "
$ws.Range("B48").Value = "'Bittensor"
$ws.Range("C48").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'405.95"
$ws.Range("E48").Value = "'  -3.78%  "
$ws.Range("E49").Value = "'  -5.88%  "
$ws.Range("D50").Value = "'142.00"
$ws.Range("E50").Value = "'  -0.04%  "
$ws.Range("E51").Value = "'  +0.13%  "
